# riska.xlsx — "Add files via upload" re-upload with refreshed repayment data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The re-uploaded copy got suffixed by Excel/Windows because a file with the
# original name already existed in the target folder.
$ws.Name = "repayment_20250916_20250916 (1)"

# Helper: write a value into a cell while forcing text storage (t="s") for
# strings that look numeric ("158,633,067.00", "1.41", ...) without leaving
# a lingering custom cell style behind.
function Set-TextCell {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Full refreshed table, row by row: A(Collector), D(Cycle), E(Repayment_amount),
# F(Pending Amount), G(Pending Amount Recovery), H(Talk_time), I(New_collections),
# J(Repayment_new_collections), K(New_collection_amount_rate), L(New_collection_count_rate)
# B (Team) and C (Cycle team code) stay "Hansyah_S2l" / "S2" for every row, unchanged.
$rows = @(
    @(2,  "Ridhoi Berkat Zebua",        3, "2,241,366.00",   "158,633,067.00", "1.41", 659, 14, 0, "0.00",  "0.00"),
    @(3,  "Adistira Winditya P",        1, "301,518.00",     "140,905,074.00", "0.21", 325, 15, 0, "0.00",  "0.00"),
    @(4,  "Yandi Nugraha",              1, "988,030.00",     "112,590,060.00", "0.88", 151, 14, 1, "10.29", "7.14"),
    @(5,  "Wasti Feronika Sihombing",   2, "375,000.00",     "145,384,324.00", "0.26", 385, 14, 0, "0.90",  "0.00"),
    @(6,  "Sucika Wardani",             3, "655,683.00",     "151,331,741.00", "0.43", 84,  14, 0, "3.17",  "0.00"),
    @(7,  "Azizah Rahmawati",           1, "1,237,616.00",   "175,179,262.00", "0.71", 65,  15, 0, "0.00",  "0.00"),
    @(8,  "Fadilah Damayanti",          1, "292,435.00",     "179,487,985.00", "0.16", 44,  14, 0, "0.00",  "0.00"),
    @(9,  "Annisa Putri Restu",         0, "0.00",           "186,099,111.00", "0.00", 449, 14, 0, "0.00",  "0.00"),
    @(10, "Riska Nurlita",              1, "220,425.00",     "188,817,452.00", "0.12", 158, 14, 0, "0.00",  "0.00"),
    @(11, "Debora Retima Sihombing",    0, "0.00",           "166,650,859.00", "0.00", 186, 14, 0, "0.00",  "0.00"),
    @(12, "Erlangga Hutama",            0, "0.00",           "113,166,920.00", "0.00", 0,   14, 0, "0.00",  "0.00"),
    @(13, "Erick Ervan Dewanggga",      0, "0.00",           "154,603,507.00", "0.00", 65,  15, 0, "0.00",  "0.00"),
    @(14, "Romli",                      0, "0.00",           "163,146,299.00", "0.00", 248, 14, 0, "0.00",  "0.00"),
    @(15, "Aldi Taufik",                1, "410,638.00",     "153,419,149.00", "0.27", 711, 15, 0, "0.00",  "0.00"),
    @(16, "Nur Halim",                  0, "0.00",           "142,201,981.00", "0.00", 134, 14, 0, "0.00",  "0.00"),
    @(17, "Axl Wicaksono",              0, "0.00",           "123,481,323.00", "0.00", 0,   14, 0, "0.00",  "0.00"),
    @(18, "Nuraini",                    0, "0.00",           "102,605,709.00", "0.00", 241, 14, 0, "0.00",  "0.00")
)

foreach ($r in $rows) {
    $row    = $r[0]
    $name   = $r[1]
    $cycle  = $r[2]
    $eVal   = $r[3]
    $fVal   = $r[4]
    $gVal   = $r[5]
    $hVal   = $r[6]
    $iVal   = $r[7]
    $jVal   = $r[8]
    $kVal   = $r[9]
    $lVal   = $r[10]

    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = "Hansyah_S2l"
    $ws.Cells.Item($row, 3).Value = "S2"
    $ws.Cells.Item($row, 4).Value = $cycle
    Set-TextCell $ws.Cells.Item($row, 5) $eVal
    Set-TextCell $ws.Cells.Item($row, 6) $fVal
    Set-TextCell $ws.Cells.Item($row, 7) $gVal
    $ws.Cells.Item($row, 8).Value  = $hVal
    $ws.Cells.Item($row, 9).Value  = $iVal
    $ws.Cells.Item($row, 10).Value = $jVal
    Set-TextCell $ws.Cells.Item($row, 11) $kVal
    Set-TextCell $ws.Cells.Item($row, 12) $lVal
}

# Selection moved from the old "F11" cursor position to the freshly
# highlighted Collector column.
$ws.Range("A2:A18").Select()
